# Generate Report for Handback
#
# The handback transform for 53b8ec75-b79d-428d-b5d6-ceae5b22d1a4 failed
# (the handback file name didn't match the handoff file name), so:
#   - its Status flips from "Ready for handoff" to "Handback transform failed"
#     (this text is a single shared string reused by the Overview sheet's
#     zh-cn/de-de status columns AND the Status column on both the zh-cn and
#     de-de detail sheets, so every one of those cells is updated here)
#   - the per-locale "Error Detail" cell for that row gets populated with the
#     mismatch explanation
#   - the "Error Detail" column is widened to fit the new text

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet, row 3 (53b8ec75...md): zh-cn / de-de status columns (E, F)
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

# Detail sheets, row 3 (53b8ec75...md): Status column (C)
$zh.Range("C3").Value = $newStatus
$de.Range("C3").Value = $newStatus

# Detail sheets, row 3: Error Detail column (P)
$zh.Range("P3").Value = "Handback file name: np2ectnh.1uj is different with handoff file name: 53b8ec75-b79d-428d-b5d6-ceae5b22d1a4.e37441b40344953a12d5fdfe95fbb492a82b8f8e.zh-cn."
$de.Range("P3").Value = "Handback file name: np2ectnh.1uj is different with handoff file name: 53b8ec75-b79d-428d-b5d6-ceae5b22d1a4.e37441b40344953a12d5fdfe95fbb492a82b8f8e.de-de."

# Widen the Error Detail column (P, the 16th column) on both detail sheets
# so the long message is readable. 39.17 is the ColumnWidth (character-unit)
# value that this workbook's other width="40" columns already use, so it
# round-trips to the raw OOXML width="40".
$zh.Columns("P").ColumnWidth = 39.17
$de.Columns("P").ColumnWidth = 39.17
